# "further cleaning to metadata"
# - The roboticS1Prep (column H) values for every data row (2-27) are
#   normalised to a single lot number "E7420" (previously most rows used
#   "E7760" while the last six rows each had their own unique lot number).
# - The formatting used for column H is refreshed to a plain Arial 10pt
#   font.
# - The active selection on the sheet is moved from column E to column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalise every roboticS1Prep (column H) cell for rows 2-27 to "E7420".
for ($row = 2; $row -le 27; $row++) {
    $ws.Range("H$row").Value = "E7420"
}

# Refresh the font used by the updated column so it renders with a plain
# Arial 10pt face (matching the rest of the sheet's body font).
$hRange = $ws.Range("H2:H27")
$hRange.Font.Name = "Arial"
$hRange.Font.Size = 10

# Move the active selection from E2:E27 to H2:H27 (active cell H2).
$ws.Range("H2:H27").Select() | Out-Null
